$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------
# The summary table's "6-cylinder" group used to be a single logical
# row spanning A9:A10 (a vertically-merged cell showing "6" once). The
# data is being re-laid out so every physical row carries its own
# cylinder value in column A and no rows are merged there anymore:
#   row 8  -> cyl=6, n=2 (engine=0)   [was old row 9]
#   row 9  -> cyl=4, n=3 (engine=1)   [was old row 8]
#   row 10 -> cyl=6, n=3 (engine=1)   [was old row 10, now gets its own "6"]
#   row 11 -> cyl=8, n=2 (unchanged)
# --------------------------------------------------------------------

# 1. Break the old A9:A10 vertical merge so each row stands alone.
$ws.Range("A9:A10").UnMerge()

# 2. A9 previously carried a unique "vertical=top" alignment (needed
#    only because it anchored a merged cell). Re-stamp A9:A10 with the
#    plain numeric-cell formatting used elsewhere in the column (copied
#    from A8) so column A looks uniform again and that now-unused style
#    stops being referenced.
$ws.Range("A8").Copy()
$ws.Range("A9:A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Write the reshuffled values, row by row.

# Row 8: cylinder 6, automatic transmission, n=2
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 110
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 2.7475
$ws.Range("G8").Value = 0.1803122292025695

# Row 9: cylinder 4, manual transmission, n=3
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 83.33333333333333
$ws.Range("E9").Value = 18.50225211517056
$ws.Range("F9").Value = 2.886666666666667
$ws.Range("G9").Value = 0.4911551010967242

# Row 10: cylinder 6, manual transmission, n=3
$ws.Range("A10").Value = 6
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 112.6666666666667
$ws.Range("E10").Value = 9.291573243177568
$ws.Range("F10").Value = 3.371666666666667
$ws.Range("G10").Value = 0.1360453355809502

# Row 11 (cylinder 8) is unchanged.

Write-Output "cylinder rows re-laid out; A9:A10 merge removed"
